{"js": "// Typo fix: the \"Ramadan is expected to end on\" date was mistyped as\n// \"Thursday 19 March 2025\" \u2014 the surrounding text already establishes\n// the 2026 Ramadan window (\"In 2026 Ramadan is expected to start on\n// Wednesday 18 February 2026\"), so the end date must read 2026 too.\nconst results = context.document.body.search(\"Thursday 19 March 2025\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"Thursday 19 March 2026\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Typo fix: the \"Ramadan is expected to end on\" date was mistyped as\n# \"Thursday 19 March 2025\" -- the surrounding text already establishes\n# the 2026 Ramadan window (\"In 2026 Ramadan is expected to start on\n# Wednesday 18 February 2026\"), so the end date must read 2026 too.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Thursday 19 March 2025\"\n$find.Replacement.Text = \"Thursday 19 March 2026\"\n$find.Forward = $true\n$find.Wrap = 0          # wdFindStop - don't wrap past the document end\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n# 1 == wdReplaceOne: replace only the (single) match found.\n$find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 1) | Out-Null\n"}
